# Fruta / hortaliza, semanal
#
# A new weekly price record for "Feria Lagunitas de Puerto Montt - Plátano"
# was inserted as row 332, pushing the former rows 332-356 down to 333-357
# (their contents are left untouched by the insert). The new row 332 gets
# its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 332; this shifts rows 332:356
# down to 333:357, preserving their existing content/formatting.
$ws.Rows(332).Insert()

# Populate the new row 332 with the new weekly record.
$ws.Range("A332").Value = 4
$ws.Range("B332").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C332").Value = "Los Lagos"
$ws.Range("D332").Value = 44578
$ws.Range("E332").Value = 10
$ws.Range("F332").Value = "Fruta"
$ws.Range("G332").Value = 100108
$ws.Range("H332").Value = "Tropicales y subtropicales"
$ws.Range("I332").Value = 100108006
$ws.Range("J332").Value = "Plátano"
$ws.Range("K332").Value = "Sin especificar"
$ws.Range("L332").Value = "Primera Pintón"
$ws.Range("M332").Value = 800
$ws.Range("N332").Value = 18000
$ws.Range("O332").Value = 19000
$ws.Range("P332").Value = 18500
$ws.Range("Q332").Value = "$/caja 20 kilos"
$ws.Range("R332").Value = "Ecuador"
$ws.Range("S332").Value = 925
$ws.Range("T332").Value = 20
